$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H11").NumberFormat = "General"
$ws.Range("H11").Value = 275

$ws.Range("H12").NumberFormat = "General"
$ws.Range("H12").Value = 425
